$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new category rows (disable duplicate categories loading
# introduced a couple of categories that used to be missing: "דוח שנתי"
# with sub-categories "תרומה מוכרת" and "ביטוח חיים").
$ws.Range("A35").Value = "דוח שנתי"
$ws.Range("B35").Value = "תרומה מוכרת"
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0

$ws.Range("A36").Value = "דוח שנתי"
$ws.Range("B36").Value = "ביטוח חיים"
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0

# Update the active selection to mirror the last edited cell, matching
# the workbook's saved cursor/view state.
[void]$ws.Range("G35").Select()
